# Apply the edit described by the diff:
#  - B5 changes from "X" to "O" (is_course_of_capacity finished)
#  - B6 changes from "X" to "O" (course_time_conflict finished)
#  - Selection moves to B8, with the view scrolled so A2 is the top-left cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "O"
$ws.Range("B6").Value = "O"

$excel.ActiveWindow.TopLeftCell = $ws.Range("A2")
$ws.Range("B8").Select()
